$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'71.453.89"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +2.34%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.649.16"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +7.91%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.01%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'587.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.89%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'180.19"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.11%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'3.638.68"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +7.76%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.623"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +4.65%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -0.02%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +1.16%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.610"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +3.17%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'49.73"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +2.24%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  -0.32%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'4.234.73"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +7.84%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'680.60"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.34%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  +4.40%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.652.74"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +7.77%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'71.576.16"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +2.58%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  +1.68%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'18.20"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +2.64%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'11.64"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +2.80%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.939"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'6.14"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +14.46%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'17.84"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +2.75%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'103.22"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +1.28%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +2.78%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  +4.48%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'10.19"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +3.50%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'35.30"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +4.84%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'9.17"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +4.50%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'7.46"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +7.33%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'4.22"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +9.72%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'579.72"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +4.14%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +2.00%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  +2.07%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'59.55"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +2.66%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'3.751.65"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +3.85%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +0.16%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +2.67%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'35.55"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +0.30%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  +3.93%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  +3.50%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  +7.82%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  +0.02%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.346"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +2.04%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'3.34"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -1.03%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'2.80"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +4.42%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.134"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +3.36%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  +3.92%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.999"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -0.24%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'134.05"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +2.57%  "
$ws.Range("E51").Style = "Normal"
